$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text-typed, matching the source data,
# so numeric-looking strings like "472.62" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.019.06"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "3.862.19"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "472.62"
$ws.Range("E5").Value = "  +10.21%  "

$ws.Range("D6").Value = "144.77"
$ws.Range("E6").Value = "  +10.16%  "

$ws.Range("E7").Value = "  +2.73%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "0.745"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -2.60%  "

$ws.Range("E11").Value = "  -8.15%  "

$ws.Range("D12").Value = "43.50"
$ws.Range("E12").Value = "  +4.44%  "

$ws.Range("D13").Value = "10.41"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "4.481.05"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "14.81"
$ws.Range("E15").Value = "  -4.46%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.867.04"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.137"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("E19").Value = "  +4.16%  "

$ws.Range("D20").Value = "67.254.74"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").Value = "432.38"
$ws.Range("E21").Value = "  +3.71%  "

$ws.Range("D22").Value = "14.95"

$ws.Range("E23").Value = "  +6.06%  "

$ws.Range("D24").Value = "88.53"
$ws.Range("E24").Value = "  +3.54%  "

$ws.Range("D25").Value = "3.59"
$ws.Range("E25").Value = "  +8.30%  "

$ws.Range("D26").Value = "38.04"
$ws.Range("E26").Value = "  +2.07%  "

$ws.Range("E27").Value = "  +5.76%  "

$ws.Range("D28").Value = "9.94"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").Value = "5.56"
$ws.Range("E29").Value = "  +3.20%  "

$ws.Range("D30").Value = "730.66"
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("D31").Value = "13.90"
$ws.Range("E31").Value = "  -2.76%  "

$ws.Range("E32").Value = "  +6.85%  "

$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("E34").Value = "  +11.00%  "

$ws.Range("E35").Value = "  +6.88%  "

$ws.Range("D36").Value = "58.31"
$ws.Range("E36").Value = "  +4.48%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").Value = "5.48"
$ws.Range("E38").Value = "  -6.52%  "

$ws.Range("E39").Value = "  +2.40%  "

$ws.Range("E40").Value = "  +7.32%  "

$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("E42").Value = "  +2.80%  "

$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0674"
$ws.Range("E43").Value = "  -7.33%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  +4.48%  "

$ws.Range("D46").Value = "3.48"
$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("E47").Value = "  +5.13%  "

$ws.Range("D48").Value = "2.16"
$ws.Range("E48").Value = "  +4.96%  "

$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("D51").Value = "143.57"
$ws.Range("E51").Value = "  +1.64%  "
